$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add "Label" header in column H, matching the style of the other headers ---
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Label"

# --- Updated D/E prediction/error values (refit results) and new H "Label" column ---
# Block 1: Iterations = 100 (rows 2-11)
$ws.Range("D2").Value  = 0.5209264890549538
$ws.Range("E2").Value  = 0.5209264890549538
$ws.Range("H2").Value  = 0

$ws.Range("D3").Value  = 0.02321475095984487
$ws.Range("E3").Value  = 0.02321475095984487
$ws.Range("H3").Value  = 0

$ws.Range("D4").Value  = 0.58965242918657
$ws.Range("E4").Value  = 0.58965242918657
$ws.Range("H4").Value  = 0

$ws.Range("D5").Value  = 0.3818128770849282
$ws.Range("E5").Value  = 0.3818128770849282
$ws.Range("H5").Value  = 0

$ws.Range("D6").Value  = 0.6124197537501553
$ws.Range("E6").Value  = 0.6124197537501553
$ws.Range("H6").Value  = 0

$ws.Range("D7").Value  = 0.6705439778232105
$ws.Range("E7").Value  = 0.3294560221767895
$ws.Range("H7").Value  = 1

$ws.Range("D8").Value  = 0.5790982383396085
$ws.Range("E8").Value  = 0.4209017616603915
$ws.Range("H8").Value  = 1

$ws.Range("D9").Value  = 0.6181890040971263
$ws.Range("E9").Value  = 0.3818109959028737
$ws.Range("H9").Value  = 1

$ws.Range("D10").Value = 0.007711851383346936
$ws.Range("E10").Value = 0.992288148616653
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.02583011479583237
$ws.Range("E11").Value = 0.9741698852041676
$ws.Range("F11").Value = 1.302705764770508
$ws.Range("H11").Value = 1

# Block 2: Iterations = 200 (rows 12-21) -- D/E unchanged here, only new H column
$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
